$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B..F shift to C..G)
$ws.Columns("B:B").Insert()

# Header for the new column: copy the header formatting (bold font,
# border, center/top alignment) from a neighboring header cell, then
# set its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B1").Value = "segments"

# Segment names (these used to live in column A) now go in column B,
# column A becomes a plain 0-based numeric index, matching the target layout.
$segments = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)

for ($i = 0; $i -lt $segments.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).ClearFormats()
    $ws.Cells.Item($row, 2).Value = $segments[$i]
}
